$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# ---------------------------------------------------------------------------
# B6 - plain text update (baseline comparison sentence appended)
# ---------------------------------------------------------------------------
$b6 = $ws.Cells.Item(6, 2)
$b6.Value = "The baseline for this target is 2015. In 2016, 91" + [char]160 + "per cent of Indigenous children were enrolled in early childhood education in the year before full time school, compared to 87 per cent in the baseline year. "

# ---------------------------------------------------------------------------
# B7 - rich text: "(ECE)" inserted in red
# ---------------------------------------------------------------------------
$b7 = $ws.Cells.Item(7, 2)
$b7.Value = "In 2016, the information available suggests that all Indigenous children were enrolled in early childhood education (ECE) in the year before full time school in Victoria, Western Australia and South Australia. Tasmania and the Australian Capital Territory met the required benchmark of 95 per cent."

$r1 = $b7.Characters(1, 116)
$r1.Font.ColorIndex = -4105
$r1.Font.Size = 12
$r1.Font.Name = "Arial"

$r2 = $b7.Characters(117, 5)
$r2.Font.Color = 255
$r2.Font.Size = 12
$r2.Font.Name = "Arial"

$r3 = $b7.Characters(122, 181)
$r3.Font.ColorIndex = -4105
$r3.Font.Size = 12
$r3.Font.Name = "Arial"

# ---------------------------------------------------------------------------
# B9 - rich text: spelled-out NECECC + a couple of red-highlighted corrections
# ---------------------------------------------------------------------------
$b9 = $ws.Cells.Item(9, 2)
$b9.Value = "Improved data quality in the 2016 National Early Childhood Education and Care Collection (NECECC) collection, resulting from revisions to the ABS data collection methodology, mean that the 2016 data are not fully comparable to the 2015 data. Siginificant changes include: amended to data linkage approach to enhance the accuracy of child counts in NECECC, and an expanded child identification strategy in the Child Care Management System (one of the source datasets) has increased the count of children enrolled in a preschool program, as all children at long day care centres (of the appropraite age) are now recorded as enrolled in a preschool program."

$s1 = $b9.Characters(1, 34)
$s1.Font.ColorIndex = -4105
$s1.Font.Size = 12
$s1.Font.Name = "Arial"

$s2 = $b9.Characters(35, 56)
$s2.Font.Color = 255
$s2.Font.Size = 12
$s2.Font.Name = "Arial"

$s3 = $b9.Characters(91, 6)
$s3.Font.ColorIndex = -4105
$s3.Font.Size = 12
$s3.Font.Name = "Arial"

$s4 = $b9.Characters(97, 2)
$s4.Font.Color = 255
$s4.Font.Size = 12
$s4.Font.Name = "Arial"

$s5 = $b9.Characters(99, 214)
$s5.Font.ColorIndex = -4105
$s5.Font.Size = 12
$s5.Font.Name = "Arial"

$s6 = $b9.Characters(313, 1)
$s6.Font.Color = 255
$s6.Font.Size = 12
$s6.Font.Name = "Arial"

$s7 = $b9.Characters(314, 341)
$s7.Font.ColorIndex = -4105
$s7.Font.Size = 12
$s7.Font.Name = "Arial"

# ---------------------------------------------------------------------------
# B10 - rich text: "ROGS" -> "Report on Government Services" (in red)
# ---------------------------------------------------------------------------
$b10 = $ws.Cells.Item(10, 2)
$b10.Value = "Preschool enrolment rates reported under the National Partnership on Universal Access to Preschool vary from the NIRA ECE rates. The NIRA and Report on Government Services enrolment rates are based on state-specific Year Before Full-Time School enrolment rates, the NIRA data also prorates Indigenous status not stated."

$t1 = $b10.Characters(1, 142)
$t1.Font.ColorIndex = -4105
$t1.Font.Size = 12
$t1.Font.Name = "Arial"

$t2 = $b10.Characters(143, 29)
$t2.Font.Color = 255
$t2.Font.Size = 12
$t2.Font.Name = "Arial"

$t3 = $b10.Characters(172, 148)
$t3.Font.ColorIndex = -4105
$t3.Font.Size = 12
$t3.Font.Name = "Arial"

Write-Host "done"
